# Add ct_loan_term_flag / ct_loan_term columns (H, I) to the "invalid" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalid")

# Header row
$ws.Range("H1").Value = "ct_loan_term_flag"
$ws.Range("I1").Value = "ct_loan_term"

# Data rows
$ws.Range("H2").Value = 900
$ws.Range("I2").Value = 0

$ws.Range("H3").Value = 900
$ws.Range("I3").Value = 1

$ws.Range("H4").Value = 900
$ws.Range("I4").Value = 1200

$ws.Range("H5").Value = 900
$ws.Range("I5").Value = 1100

$ws.Range("H6").Value = 900
$ws.Range("I6").Value = 0.5

$ws.Range("I7").Value = 1

$ws.Range("H8").Value = 988
$ws.Range("I8").Value = 1

$ws.Range("H9").Value = 999
$ws.Range("I9").Value = 36

$ws.Range("H10").Value = 988

$ws.Range("H11").Value = 999

# Reflect the final selection left behind in the author's editing session.
$ws.Activate() | Out-Null
$ws.Range("I15").Select() | Out-Null
